$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.689.83"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.597.85"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'0.511"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.0618"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").Value = "'19.73"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'0.0838"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.821.47"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.590.90"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "'65.15"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "26.672.95"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "'209.88"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").Value = "'8.93"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'146.67"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'7.17"
$ws.Range("E27").Value = "  -4.20%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "'15.31"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'0.666"
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "1.297.30"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D38").Value = "'0.0171"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'63.80"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "1.734.22"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").Value = "'0.886"
$ws.Range("E46").Value = "  +9.85%  "
$ws.Range("D47").Value = "'90.24"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'1.62"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "'0.101"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  +0.39%  "
